$wb = $excel.ActiveWorkbook

# Delete the obsolete "Desarquivamentos Pendentes" sheet entirely.
$old = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$old.Delete()

# Rename "Paineis DARQ" -> "PAINEIS DARQ" (keep it the active/selected tab)
$paineis = $wb.Worksheets.Item("Paineis DARQ")
$paineis.Name = "PAINEIS DARQ"
$paineis.Activate()

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$recolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$recolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"
